$d = $word.ActiveDocument

# --- Locate the target paragraph: the one ending in "...是jsonl。" ---
$targetRange = $d.Content
$found = $targetRange.Find.Execute("错误的格式是指我之前用的是json数据的格式", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target paragraph"
}
$null = $targetRange.Expand(4)
$paraStart = $targetRange.Start
$paraEnd = $targetRange.End

# --- Build the replacement paragraph-end fragment: new run + moved bookmark ---
$newRunText = "Json对象中system和tools可有可无，如果要添加这两个属性，在data_info中columns中添加这两个属性即可。"

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="53F47BBE"><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>' + $newRunText + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Insert the fragment as a brand-new paragraph right after the target paragraph
# (position: just before the paragraph's own trailing pilcrow).
$insertionPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
$null = $insertionPoint.InsertXML($fragment)

# Merge the newly inserted paragraph back into the target paragraph by deleting
# the paragraph mark that separates them. Word transplants the deleted pilcrow's
# *successor* paragraph mark (ours, carrying pPr/rPr + trailing bookmark) onto
# the merged paragraph -- exactly matching the diff's target shape.
$paraMark = $d.Range($paraEnd - 1, $paraEnd)
$paraMark.Delete()

# --- Remove the old _GoBack bookmark that used to sit at the end of the document ---
# There are now two "_GoBack" bookmarks: the one we just inserted (inside the
# merged paragraph) and the original one at the end of the document. Walk the
# document's raw bookmark list (by position) and delete the one that is NOT
# inside our merged paragraph.
$mergedParaEnd = $paraEnd + (Get-StringLength $newRunText)
$goBackStart = $null
try {
    $probe = $d.Bookmarks.Item("_GoBack")
    $goBackStart = $probe.Start
} catch {
}
if ($goBackStart -ne $null -and $goBackStart -gt $mergedParaEnd) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
